# Auto-generated Excel COM-interop script applying the meteocat daily-summary refresh
# (DATA_EXTRACCIO timestamps + a handful of recomputed daily stats) for run at 2026-02-24 06:50.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell holding the plain (unformatted, style index 3) "General" style used to restore
# formatting on cells whose new text looks like a percentage (Excel would otherwise silently
# reinterpret "40%" as the number 0.4 with a Percent number format).
$formatDonor = $ws.Range("C2")

# --- Plain text/number-like values (dates, °C, hPa, km/h, etc.) ---------------------------
$ws.Range("E2").Value = '2026-02-24 06:48:24'
$ws.Range("E3").Value = '2026-02-24 06:48:27'
$ws.Range("O3").Value = '2.4 °C'
$ws.Range("E4").Value = '2026-02-24 06:48:29'
$ws.Range("N4").Value = '4.7 °C 6:24 TU'
$ws.Range("O4").Value = '6.7 °C'
$ws.Range("E5").Value = '2026-02-24 06:48:32'
$ws.Range("E6").Value = '2026-02-24 06:48:34'
$ws.Range("O6").Value = '9.1 °C'
$ws.Range("E7").Value = '2026-02-24 06:48:37'
$ws.Range("N7").Value = '10.8 °C 6:29 TU'
$ws.Range("O7").Value = '12.0 °C'
$ws.Range("E8").Value = '2026-02-24 06:48:39'
$ws.Range("N8").Value = '13.3 °C 6:29 TU'
$ws.Range("O8").Value = '14.8 °C'
$ws.Range("E9").Value = '2026-02-24 06:48:42'
$ws.Range("O9").Value = '5.4 °C'
$ws.Range("E10").Value = '2026-02-24 06:48:44'
$ws.Range("E11").Value = '2026-02-24 06:48:47'
$ws.Range("N11").Value = '1.4 °C 6:04 TU'
$ws.Range("O11").Value = '2.5 °C'
$ws.Range("E12").Value = '2026-02-24 06:48:49'
$ws.Range("O12").Value = '5.6 °C'
$ws.Range("E13").Value = '2026-02-24 06:48:51'
$ws.Range("O13").Value = '-1.7 °C'
$ws.Range("E14").Value = '2026-02-24 06:48:54'
$ws.Range("E15").Value = '2026-02-24 06:48:56'
$ws.Range("O15").Value = '5.6 °C'
$ws.Range("E16").Value = '2026-02-24 06:48:59'
$ws.Range("E17").Value = '2026-02-24 06:49:01'
$ws.Range("E18").Value = '2026-02-24 06:49:04'
$ws.Range("J18").Value = '1022.5 hPa'
$ws.Range("N18").Value = '1.2 °C 6:21 TU'
$ws.Range("O18").Value = '2.8 °C'
$ws.Range("E19").Value = '2026-02-24 06:49:06'
$ws.Range("E20").Value = '2026-02-24 06:49:08'
$ws.Range("N20").Value = '-0.8 °C 6:22 TU'
$ws.Range("E21").Value = '2026-02-24 06:49:11'
$ws.Range("N21").Value = '1.0 °C 6:02 TU'
$ws.Range("O21").Value = '3.2 °C'
$ws.Range("E22").Value = '2026-02-24 06:49:14'
$ws.Range("N22").Value = '1.5 °C 6:24 TU'
$ws.Range("O22").Value = '2.7 °C'
$ws.Range("E23").Value = '2026-02-24 06:49:16'
$ws.Range("O23").Value = '3.8 °C'
$ws.Range("E24").Value = '2026-02-24 06:49:19'
$ws.Range("N24").Value = '0.7 °C 6:29 TU'
$ws.Range("O24").Value = '2.9 °C'
$ws.Range("E25").Value = '2026-02-24 06:49:21'
$ws.Range("O25").Value = '4.9 °C'
$ws.Range("E26").Value = '2026-02-24 06:49:23'
$ws.Range("E27").Value = '2026-02-24 06:49:26'
$ws.Range("N27").Value = '3.4 °C 6:04 TU'
$ws.Range("E28").Value = '2026-02-24 06:49:28'
$ws.Range("J28").Value = '1023.7 hPa'
$ws.Range("O28").Value = '3.6 °C'
$ws.Range("E29").Value = '2026-02-24 06:49:31'
$ws.Range("N29").Value = '3.0 °C 6:01 TU'
$ws.Range("O29").Value = '4.8 °C'
$ws.Range("E30").Value = '2026-02-24 06:49:33'
$ws.Range("N30").Value = '7.4 °C 6:14 TU'
$ws.Range("O30").Value = '8.9 °C'
$ws.Range("E31").Value = '2026-02-24 06:49:36'
$ws.Range("O31").Value = '14.7 °C'
$ws.Range("E32").Value = '2026-02-24 06:49:39'
$ws.Range("E33").Value = '2026-02-24 06:49:41'
$ws.Range("O33").Value = '1.9 °C'
$ws.Range("E34").Value = '2026-02-24 06:49:44'
$ws.Range("O34").Value = '2.7 °C'
$ws.Range("E35").Value = '2026-02-24 06:49:47'
$ws.Range("E36").Value = '2026-02-24 06:49:49'
$ws.Range("N36").Value = '5.7 °C 6:25 TU'
$ws.Range("O36").Value = '7.9 °C'
$ws.Range("E37").Value = '2026-02-24 06:49:52'
$ws.Range("J37").Value = '1027.3 hPa'
$ws.Range("L37").Value = '10.8 km/h - 212º 6:28 TU'
$ws.Range("E38").Value = '2026-02-24 06:49:55'
$ws.Range("N38").Value = '4.2 °C 6:01 TU'
$ws.Range("O38").Value = '6.4 °C'
$ws.Range("E39").Value = '2026-02-24 06:49:57'
$ws.Range("N39").Value = '2.4 °C 6:29 TU'
$ws.Range("E40").Value = '2026-02-24 06:50:00'
$ws.Range("N40").Value = '-0.4 °C 6:29 TU'
$ws.Range("O40").Value = '1.0 °C'
$ws.Range("E41").Value = '2026-02-24 06:50:02'
$ws.Range("J41").Value = '1022.1 hPa'
$ws.Range("N41").Value = '4.3 °C 6:21 TU'
$ws.Range("O41").Value = '6.9 °C'
$ws.Range("E42").Value = '2026-02-24 06:50:05'
$ws.Range("O42").Value = '6.6 °C'
$ws.Range("E43").Value = '2026-02-24 06:50:07'
$ws.Range("N43").Value = '2.0 °C 6:18 TU'
$ws.Range("O43").Value = '3.9 °C'
$ws.Range("E44").Value = '2026-02-24 06:50:10'
$ws.Range("E45").Value = '2026-02-24 06:50:13'
$ws.Range("J45").Value = '1025.4 hPa'
$ws.Range("L45").Value = '16.6 km/h - 138º 6:05 TU'
$ws.Range("O45").Value = '4.6 °C'
$ws.Range("E46").Value = '2026-02-24 06:50:15'
$ws.Range("N46").Value = '0.5 °C 6:12 TU'
$ws.Range("O46").Value = '2.1 °C'

# --- Percentage-looking text values -------------------------------------------------------
# Assigned with the cell pre-formatted as Text so Excel stores the literal string instead of
# converting it to a numeric percentage, then the format is restored to the original plain
# style via a copy/paste-special of just the formatting (so the stored value stays untouched).
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '40%'
$formatDonor.Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null

$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '33%'
$formatDonor.Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null

$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '93%'
$formatDonor.Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null

$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '92%'
$formatDonor.Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null

$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '17%'
$formatDonor.Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null

$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '42%'
$formatDonor.Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4122) | Out-Null

$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '26%'
$formatDonor.Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null

$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '52%'
$formatDonor.Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null

$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '85%'
$formatDonor.Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null

$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '80%'
$formatDonor.Copy() | Out-Null
$ws.Range("H41").PasteSpecial(-4122) | Out-Null

$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '51%'
$formatDonor.Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4122) | Out-Null

$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '59%'
$formatDonor.Copy() | Out-Null
$ws.Range("H45").PasteSpecial(-4122) | Out-Null

